$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.08204471844477863
$ws.Range("E2").Value = 0.08204471844477863

# Row 3
$ws.Range("D3").Value = 0.2798291202851274
$ws.Range("E3").Value = 0.2798291202851274

# Row 4
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = 0.01694158468343083
$ws.Range("E4").Value = 0.01694158468343083

# Row 5
$ws.Range("D5").Value = 0.2332477623383898
$ws.Range("E5").Value = 0.2332477623383898

# Row 6
$ws.Range("D6").Value = 0.1996456718942284
$ws.Range("E6").Value = 0.1996456718942284

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.01431871118954079
$ws.Range("E7").Value = 0.9856812888104592

# Row 8
$ws.Range("D8").Value = 0.0007832998432050463
$ws.Range("E8").Value = 0.999216700156795

# Row 9
$ws.Range("D9").Value = 0.03516121175312053
$ws.Range("E9").Value = 0.9648387882468795

# Row 10
$ws.Range("D10").Value = 0.05588387330400575
$ws.Range("E10").Value = 0.9441161266959942
$ws.Range("F10").Value = 2.061080694198608

# Row 11
$ws.Range("D11").Value = 0.04947099547241984
$ws.Range("E11").Value = 0.04947099547241984

# Row 12
$ws.Range("D12").Value = 0.2613839120392689
$ws.Range("E12").Value = 0.2613839120392689

# Row 13
$ws.Range("C13").Value = $true
$ws.Range("D13").Value = 0.02236916021366518
$ws.Range("E13").Value = 0.02236916021366518

# Row 14
$ws.Range("D14").Value = 0.1918938443997983
$ws.Range("E14").Value = 0.1918938443997983

# Row 15
$ws.Range("D15").Value = 0.1475817249778648
$ws.Range("E15").Value = 0.1475817249778648

# Row 16
$ws.Range("D16").Value = 0.007898591266148804
$ws.Range("E16").Value = 0.9921014087338512

# Row 17
$ws.Range("D17").Value = 0.02619233908752912
$ws.Range("E17").Value = 0.9738076609124708

# Row 18
$ws.Range("D18").Value = 0.07495124491579559
$ws.Range("E18").Value = 0.9250487550842044

# Row 19
$ws.Range("D19").Value = 0.164004544067148
$ws.Range("E19").Value = 0.8359954559328521
$ws.Range("F19").Value = 1.514579296112061
$ws.Range("G19").Value = 0.5555555555555556
